# Daily auto push: insert one new data row for 2026/01/22 (木) into the
# time-series table on Sheet1. The new row is inserted immediately before
# the existing row 676 ("2026/12/29"), which pushes every row from the old
# 676 through 717 down by one (now 677-718), and extends the used range
# from A1:D717 to A1:D718.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 676..717 down to 677..718 by inserting a blank row at 676.
$ws.Rows("676:676").Insert()

# Populate the newly inserted row 676 with the new observation:
# date 2026/01/22, weekday 木, time 13, ranking 16.
#
# Columns A/B hold literal text (not real Excel dates) in every row of
# this table, e.g. rows 673-675 already contain the same "2026/01/22" /
# "木" strings. Assigning a date-shaped string straight to .Value would
# make Excel auto-convert it to a date serial (and tag the cell with a
# date number format), which would NOT match the plain text the rest of
# the column uses. Copying the existing matching cells instead keeps the
# new cells as plain text, identical in kind to their neighbours, with no
# incidental number-format/style changes.
$ws.Range("A673").Copy($ws.Range("A676"))
$ws.Range("B673").Copy($ws.Range("B676"))

# Columns C/D are genuine numbers - plain assignment is fine for these.
$ws.Range("C676").Value = 13
$ws.Range("D676").Value = 16
